$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.980.49"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "2.234.74"
$ws.Range("E3").Value = "  -3.35%  "
$ws.Range("D5").Value = "'292.87"
$ws.Range("E5").Value = "  -4.70%  "
$ws.Range("D6").Value = "'86.43"
$ws.Range("E6").Value = "  +5.83%  "
$ws.Range("D7").Value = "'0.513"
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.470"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "'30.74"
$ws.Range("E10").Value = "  +7.16%  "
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").Value = "'47.24"
$ws.Range("E12").Value = "  -9.67%  "
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").Value = "'6.39"
$ws.Range("E14").Value = "  +2.48%  "
$ws.Range("D15").Value = "2.576.71"
$ws.Range("E15").Value = "  -3.64%  "
$ws.Range("D16").Value = "'14.17"
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("D17").Value = "2.214.54"
$ws.Range("E17").Value = "  -4.67%  "
$ws.Range("D18").Value = "'0.726"
$ws.Range("E18").Value = "  -1.62%  "
$ws.Range("D19").Value = "39.884.09"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("D20").Value = "0.0₃0894"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").Value = "'5.80"
$ws.Range("E21").Value = "  -2.14%  "
$ws.Range("D22").Value = "'10.68"
$ws.Range("E22").Value = "  +4.42%  "
$ws.Range("D23").Value = "'65.48"
$ws.Range("E23").Value = "  -2.70%  "
$ws.Range("D24").Value = "'234.90"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "'2.43"
$ws.Range("E26").Value = "  -2.40%  "
$ws.Range("D27").Value = "'1.84"
$ws.Range("E27").Value = "  +3.88%  "
$ws.Range("D28").Value = "'22.94"
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("E29").Value = "  +1.27%  "
$ws.Range("D30").Value = "'9.25"
$ws.Range("E30").Value = "  +2.44%  "
$ws.Range("D31").Value = "'33.43"
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("D32").Value = "'154.60"
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").Value = "'4.86"
$ws.Range("E34").Value = "  -1.74%  "
$ws.Range("D35").Value = "'0.0712"
$ws.Range("E35").Value = "  +1.69%  "
$ws.Range("E36").Value = "  -2.63%  "
$ws.Range("D37").Value = "'16.58"
$ws.Range("E37").Value = "  +9.77%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").Value = "'0.0999"
$ws.Range("E39").Value = "  +3.53%  "
$ws.Range("D40").Value = "'2.69"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").Value = "'1.67"
$ws.Range("E41").Value = "  +2.31%  "
$ws.Range("D42").Value = "'3.78"
$ws.Range("E42").Value = "  +2.58%  "
$ws.Range("D43").Value = "1.956.19"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("E44").Value = "  -2.93%  "
$ws.Range("E45").Value = "  +5.08%  "
$ws.Range("D46").Value = "'9.58"
$ws.Range("E46").Value = "  +4.08%  "
$ws.Range("D47").Value = "'16.22"
$ws.Range("E47").Value = "  -2.91%  "
$ws.Range("D48").Value = "'2.60"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").Value = "2.450.71"
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").Value = "'71.01"
$ws.Range("E50").Value = "  +4.01%  "
$ws.Range("D51").Value = "'1.46"
$ws.Range("E51").Value = "  +10.07%  "
